$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells removed in target ---
$ws.Range("E2").ClearContents()
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# --- Set changed / added cell values ---
$ws.Range("D2").Value = -0.108
$ws.Range("F2").Value = -0.07530000000000001
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -1513.6
$ws.Range("L2").Value = -0.425228262396404
$ws.Range("M2").Value = 27
$ws.Range("N2").Value = 0.00259927797833935
$ws.Range("O2").Value = -0.01783826638477801
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 0.00259927797833935
$ws.Range("R2").Value = -0.01783826638477801
$ws.Range("U2").Value = 30165.5
$ws.Range("V2").Value = 2.904019253910951
$ws.Range("W2").Value = -0.05194184525606198
$ws.Range("X2").Value = 0.08654176508662438
$ws.Range("Y2").Value = -0.1384836103426864
$ws.Range("Z2").Value = 0.09070315008383577
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.05143221905562387
$ws.Range("AC2").Value = -0.05143221905562387
$ws.Range("AD2").Value = 26073.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 26073.6
$ws.Range("AG2").Value = -4091.900000000001
$ws.Range("AH2").Value = 0.7151073335692011
$ws.Range("AI2").Value = 0.4743768182388804
$ws.Range("AJ2").Value = -0.649961878137112
$ws.Range("AK2").Value = -0.1650066133298923
$ws.Range("B3").Value = "AIB Group plc (ISE:A5G)"
$ws.Range("D3").Value = -0.194
$ws.Range("F3").Value = -0.0595
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -826
$ws.Range("L3").Value = -0.5933481790101286
$ws.Range("M3").Value = 14.6
$ws.Range("N3").Value = 0.002615784287377945
$ws.Range("O3").Value = -0.01767554479418886
$ws.Range("P3").Value = 14.6
$ws.Range("Q3").Value = 0.002615784287377945
$ws.Range("R3").Value = -0.01767554479418886
$ws.Range("U3").Value = 17777.4
$ws.Range("V3").Value = 3.185057780166622
$ws.Range("W3").Value = -0.05194184525606198
$ws.Range("X3").Value = 0.08356422632073331
$ws.Range("Y3").Value = -0.1355060715767953
$ws.Range("Z3").Value = 0.08292924118046538
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04966946306528673
$ws.Range("AC3").Value = -0.04966946306528673
$ws.Range("AD3").Value = 11457.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 11457.3
$ws.Range("AG3").Value = -6320.100000000002
$ws.Range("AH3").Value = 0.6724241143742516
$ws.Range("AI3").Value = 0.4250402510776901
$ws.Range("AJ3").Value = 8.556864337936615
$ws.Range("AK3").Value = -0.688584066939772
$ws.Range("B4").Value = "Bank of Ireland Group plc (ISE:BIRG)"
$ws.Range("D4").Value = -0.108
$ws.Range("F4").Value = -0.07530000000000001
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = -637
$ws.Range("L4").Value = -0.3498462214411248
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 10561.6
$ws.Range("V4").Value = 2.441085378819396
$ws.Range("W4").Value = -0.06047257848619193
$ws.Range("X4").Value = 0.1087529485405322
$ws.Range("Y4").Value = -0.1692255270267241
$ws.Range("Z4").Value = 0.08850693161711809
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.05143221905562387
$ws.Range("AC4").Value = -0.05143221905562387
$ws.Range("AD4").Value = 13570.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 13570.8
$ws.Range("AG4").Value = 3009.199999999999
$ws.Range("AH4").Value = 0.7582553890509235
$ws.Range("AI4").Value = 0.5474657501089218
$ws.Range("AJ4").Value = 0.4102074756672754
$ws.Range("AK4").Value = 0.2115162931931284
$ws.Range("F5").Value = -0.156
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -50.6
$ws.Range("L5").Value = -0.1459896133871898
$ws.Range("M5").Value = 12.4
$ws.Range("N5").Value = 0.02586566541510221
$ws.Range("O5").Value = -0.2450592885375494
$ws.Range("P5").Value = 12.4
$ws.Range("Q5").Value = 0.02586566541510221
$ws.Range("R5").Value = -0.2450592885375494
$ws.Range("U5").Value = 1826.5
$ws.Range("V5").Value = 3.80997079682937
$ws.Range("W5").Value = -0.02236760675448678
$ws.Range("X5").Value = 0.08654176508662438
$ws.Range("Y5").Value = -0.1089093718411112
$ws.Range("Z5").Value = 0.1839312247930376
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.05770162518161664
$ws.Range("AC5").Value = -0.05770162518161664
$ws.Range("AD5").Value = 1045.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 1045.5
$ws.Range("AG5").Value = -781
$ws.Range("AH5").Value = 0.68561872909699
$ws.Range("AI5").Value = 0.3247196943814641
$ws.Range("AJ5").Value = 2.589522546419098
$ws.Range("AK5").Value = -0.560579959804766
$ws.Range("D5").Value = -0.0883
